# Sample_import_banking_file_en.xlsx — strip the Arabic translations from
# the bilingual (EN/AR) header row, leaving English-only header labels.
# Visible cell positions (A1:E1) are unchanged; only the trailing
# "/<Arabic>" portion of each label is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "* Date (dd/mm/YYYY)"
$ws.Range("B1").Value = "* Withdrawals"
$ws.Range("C1").Value = "* Deposits"
$ws.Range("D1").Value = " * Payee"
$ws.Range("E1").Value = "Description"
